$wb = $excel.ActiveWorkbook

# --- Update the daily conversion text on "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.17 = 16570.88 pesos`n✅ 16570.88 pesos = 4.15 = 912.83 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 240
$wsTasas.Range("O10").Value = 3977.01

$wsTasas.Range("N12").Value = 3990
$wsTasas.Range("O12").Value = 219.795
